$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old SUM/SQRT demonstration block that lived in columns E:G (rows 17-20)
$ws.Range("E17:G20").ClearContents()

# New parts added to the BOM (rows 15-18)
$ws.Range("A15").Value = "RX"
$ws.Range("A16").Value = "N channel MOSFET"
$ws.Range("A17").Value = "p CHANNEL MOSFET"
$ws.Range("B17").Value = "https://www.digikey.com/product-detail/en/infineon-technologies/IRF9540NSTRLPBF/IRF9540NSTRLPBFTR-ND/1928217"
$ws.Range("B16").Value = "https://www.digikey.com/product-detail/en/toshiba-semiconductor-and-storage/T2N7002AK,LM/T2N7002AKLMCT-ND/5298039"
$ws.Range("A18").Value = "Button"
$ws.Range("B18").Value = "https://www.digikey.com/product-detail/en/e-switch/KS-01Q-02/EG4792-ND/2116271"

# Move the selection, matching the saved view state in the file
$ws.Range("C22").Select()
